# Update the IRR/cash-flow figures in column J (rows 4-38) on Sheet1, then
# leave the active selection on Q16 (matching the author's last-saved cursor
# position) as recorded in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("J4").Value = 2501.9953071004629
$ws.Range("J5").Value = 3835.6244504194347
$ws.Range("J6").Value = 3285.1223048517027
$ws.Range("J7").Value = 3937.1846059278014
$ws.Range("J8").Value = 4170.6304024117471
$ws.Range("J9").Value = 4248.1864762641917
$ws.Range("J10").Value = 4454.0089124215929
$ws.Range("J11").Value = 4278.6716997757094
$ws.Range("J12").Value = 4259.000194904238
$ws.Range("J13").Value = 4383.4562577475817
$ws.Range("J14").Value = 3832.9880163038811
$ws.Range("J15").Value = 3913.2565172921413
$ws.Range("J16").Value = 3937.7147867173735
$ws.Range("J17").Value = 2101.0446766878317
$ws.Range("J18").Value = 2935.3275547712306
$ws.Range("J19").Value = 2176.8108881056482
$ws.Range("J20").Value = 2023.8412444417288
$ws.Range("J21").Value = 1870.1047578422101
$ws.Range("J22").Value = 1961.2355100722823
$ws.Range("J23").Value = 807.89927763605942
$ws.Range("J24").Value = 854.24973682712641
$ws.Range("J25").Value = 873.16574575509287
$ws.Range("J26").Value = 790.38698278973141
$ws.Range("J27").Value = -116.07390559453955
$ws.Range("J28").Value = -227.95215347026419
$ws.Range("J29").Value = -184.46003839554967
$ws.Range("J30").Value = -106.31302552011442
$ws.Range("J31").Value = 13.935187786349861
$ws.Range("J32").Value = 81.213714187589304
$ws.Range("J33").Value = -210.5734693421868
$ws.Range("J34").Value = -583.55332889163435
$ws.Range("J35").Value = -978.1043834942775
$ws.Range("J36").Value = -931.82308742781004
$ws.Range("J37").Value = -896.1880222748066
$ws.Range("J38").Value = -248.57038566350593

# Move the selection to Q16, matching the saved cursor position in the diff.
$ws.Range("Q16").Select()
